# Insert a new data row at row 24, shifting existing rows 24-55 down to 25-56,
# and populate the new row 24 with a new "Arveja Verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 24 (pushes old rows 24.. down by one)
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record's values
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 4).Value = 44581
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = 100112022
$ws.Cells.Item(24, 7).Value = "Arveja Verde"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 80
$ws.Cells.Item(24, 11).Value = 24000
$ws.Cells.Item(24, 12).Value = 25000
$ws.Cells.Item(24, 13).Value = 24500
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(24, 16).Value = 980
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
